$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the "plus" header (AA1) to "plus/minus"
$ws.Range("AA1").Value = "plus/minus"

# New team rows appended to the table (2018-2019 NBA season summary)
$teams = @(
    @{ Row = 24; Name = "Denver Nuggets";      Vals = @(82,54,28,0.659,48.1,110.7,41.9,90.0,46.6,11.0,31.4,35.1,15.8,20.9,75.5,11.9,34.5,46.4,27.4,13.4,7.7,4.4,5.0,20.0,20.4,4.0) },
    @{ Row = 25; Name = "Dallas Mavericks";    Vals = @(82,33,49,0.402,48.2,108.9,38.8,86.9,44.7,12.5,36.6,34.0,18.8,25.3,74.2,10.1,35.2,45.3,23.4,14.2,6.5,4.3,4.5,20.1,23.2,-1.3) },
    @{ Row = 26; Name = "Cleveland Cavaliers"; Vals = @(82,19,63,0.232,48.2,104.5,38.9,87.6,44.4,10.3,29.1,35.5,16.4,20.7,79.2,10.7,31.9,42.7,20.7,13.5,6.5,2.4,5.6,20.0,19.4,-9.6) },
    @{ Row = 27; Name = "Chicago Bulls";       Vals = @(82,22,60,0.268,48.5,104.9,39.8,87.9,45.3,9.1,25.9,35.1,16.2,20.7,78.3,8.8,34.1,42.9,21.9,14.1,7.4,4.3,5.8,20.3,18.7,-8.4) },
    @{ Row = 28; Name = "Charlotte Hornets";   Vals = @(82,39,43,0.476,48.4,110.7,40.2,89.8,44.8,11.9,33.9,35.1,18.4,23.1,79.7,9.9,33.9,43.8,23.2,12.2,7.2,4.9,6.0,18.9,20.6,-1.1) },
    @{ Row = 29; Name = "Brooklyn Nets";       Vals = @(82,42,40,0.512,48.7,112.2,40.3,89.7,44.9,12.8,36.2,35.3,19.0,25.5,74.5,11.0,35.6,46.6,23.8,15.1,6.6,4.1,5.3,21.5,22.0,-0.1) },
    @{ Row = 30; Name = "Boston Celtics";      Vals = @(82,49,33,0.598,48.2,112.4,42.1,90.5,46.5,12.6,34.5,36.5,15.6,19.5,80.2,9.8,34.7,44.5,26.3,12.8,8.6,5.3,3.9,20.4,19.5,4.4) },
    @{ Row = 31; Name = "Atlanta Hawks";       Vals = @(82,29,53,0.354,48.4,113.3,41.4,91.8,45.1,13.0,37.0,35.2,17.6,23.4,75.2,11.6,34.5,46.1,25.8,17.0,8.2,5.1,5.5,23.6,22.2,-6.0) }
)

$cols = @("B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V","W","X","Y","Z","AA")

foreach ($team in $teams) {
    $r = $team.Row
    $ws.Range("A$r").Value = $team.Name
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Range("$($cols[$i])$r").Value = $team.Vals[$i]
    }
}
